# Change SRS front page
#   1. Remove the "Kshitij Nigam" paragraph from the title page.
#   2. Re-layout the Table-of-Contents table: narrow col 1, widen col 2,
#      fill in the S.No./Title/Page No. header and the TOC entries, and
#      append a 6th data row.

$d = $word.ActiveDocument

# --- 1. Delete the whole paragraph that only contains "Kshitij Nigam" ---
$find = $d.Content.Find
$find.ClearFormatting()
if ($find.Execute("Kshitij Nigam")) {
    $rng = $find.Parent
    $rng.Expand(4) | Out-Null  # wdParagraph - grab the paragraph incl. its mark
    $rng.Delete()
}

# --- 2. Table of contents table -----------------------------------------
$t = $d.Tables.Item(1)

# Resize the grid: col1 narrower, col2 wider, col3 untouched
$t.Columns.Item(1).Width = 65.15   # 1303 twips
$t.Columns.Item(2).Width = 256.1   # 5122 twips

# Header row
$t.Cell(1, 2).Range.Text = "Title"
$t.Cell(1, 3).Range.Text = "Page No."

# Row 2
$t.Cell(2, 1).Range.Text = "1."
$t.Cell(2, 2).Range.Text = "Introduction"

# Row 3
$t.Cell(3, 1).Range.Text = "2."
$t.Cell(3, 2).Range.Text = "Overall description"

# Row 4
$t.Cell(4, 1).Range.Text = "3."
$t.Cell(4, 2).Range.Text = "Specific Requirments"

# Row 5
$t.Cell(5, 1).Range.Text = "4."

# Row 6
$t.Cell(6, 1).Range.Text = "5."

# Append a new 7th row
$t.Rows.Add() | Out-Null
$t.Cell(7, 1).Range.Text = "6."

Write-Host "SRS front page updated"
